$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells (Wins, Losses, Ties) right after the
# existing "Unnamed: 28" column (AC), matching the style of the other
# header cells (bold/centered/bordered, same format as A1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill season-record values for every player row (2-50): 94 wins,
# 67 losses, 0 ties for the 2016 season.
$ws.Range("AD2:AD50").Value = 94
$ws.Range("AE2:AE50").Value = 67
$ws.Range("AF2:AF50").Value = 0
